$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The module list shrank down to a single row (FASB / Yes) and the helper
# "Yes"/"No" list that feeds the ExecutionFlag data-validation dropdown moved
# up from rows 49:50 to rows 43:44 (it always sits right after the last used
# row). Remove the now-obsolete module rows (3-8) - Excel shifts the C49:C50
# helper cells up automatically when the rows above them disappear.
$ws.Range("A3:B8").EntireRow.Delete() | Out-Null

# Row 2 becomes the single remaining module entry.
$ws.Range("A2").Value = "Yes"
$ws.Range("B2").Value = "FASB"

# B2 previously relied on the row-1 bottom border for its top edge; make it
# match A2's own all-round thin border explicitly (copy A2's formatting onto
# B2) instead of leaving it on the old border-without-top style.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Re-point the ExecutionFlag dropdown's list source at the relocated C43:C44
# cells, and restore the validation to cover the whole column (deleting rows
# shrank the sqref's lower bound, so drop and re-add it).
$ws.Range("A1:A1048576").Validation.Delete() | Out-Null
$ws.Range("A1:A1048576").Validation.Add(3, 1, 1, "=`$C`$43:`$C`$44") | Out-Null

# Match the author's final selection (whole of row 2 highlighted).
$ws.Rows("2:2").Select() | Out-Null
